$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new station/platform rows 224-323 ---

# Row 224: Nar Nar Goon
$ws.Range("A224").Formula = '=TRANSPOSE({"Nar Nar Goon","Tynong","Garfield","Bunyip","Longwarry","Drouin","Warragul","Yarragon","Trafalgar","Moe","Morwell","Traralgon","Rosedale","Sale","Stratford","Bairnsdale"})'
$ws.Range("B224").Value = 'L'
$ws.Range("C224").Value = 'R'

# Row 225: Tynong
$ws.Range("B225").Value = 'L'
$ws.Range("C225").Value = 'R'

# Row 226: Garfield
$ws.Range("B226").Value = 'R'
$ws.Range("C226").Value = 'L'

# Row 227: Bunyip
$ws.Range("B227").Value = 'R'
$ws.Range("C227").Value = 'L'

# Row 228: Longwarry
$ws.Range("B228").Value = 'R'
$ws.Range("C228").Value = 'L'

# Row 229: Drouin
$ws.Range("B229").Value = 'L'
$ws.Range("C229").Value = 'R'

# Row 230: Warragul
$ws.Range("B230").Value = 'L'
$ws.Range("C230").Value = 'R'

# Row 231: Yarragon
$ws.Range("B231").Value = 'R'
$ws.Range("C231").Value = 'L'

# Row 232: Trafalgar
$ws.Range("B232").Value = 'R'
$ws.Range("C232").Value = 'L'

# Row 233: Moe
$ws.Range("B233").Value = 'R'

# Row 234: Morwell
$ws.Range("B234").Value = 'L'

# Row 235: Traralgon
$ws.Range("B235").Value = 'L'
$ws.Range("C235").Value = 'R'

# Row 236: Rosedale
$ws.Range("B236").Value = 'R'

# Row 237: Sale
$ws.Range("B237").Value = 'R'

# Row 238: Stratford
$ws.Range("B238").Value = 'L'

# Row 239: Bairnsdale
$ws.Range("B239").Value = 'L'

# Row 240: Donnybrook
$ws.Range("A240").Formula = '=TRANSPOSE({"Donnybrook","Wallan","Heathcote Junction","Wandong","Kilmore East","Broadford","Tallarook","Seymour","Avenel","Euroa","Violet Town","Benalla","Wangaratta","Springhurst","Chiltern","Wodonga","Albury"})'
$ws.Range("B240").Value = 'R'
$ws.Range("C240").Value = 'L'

# Row 241: Wallan
$ws.Range("B241").Value = 'R'
$ws.Range("C241").Value = 'L'

# Row 242: Heathcote Junction
$ws.Range("B242").Value = 'R'
$ws.Range("C242").Value = 'L'

# Row 243: Wandong
$ws.Range("B243").Value = 'R'
$ws.Range("C243").Value = 'L'

# Row 244: Kilmore East
$ws.Range("B244").Value = 'R'
$ws.Range("C244").Value = 'L'

# Row 245: Broadford
$ws.Range("B245").Value = 'R'
$ws.Range("C245").Value = 'L'

# Row 246: Tallarook
$ws.Range("B246").Value = 'R'
$ws.Range("C246").Value = 'L'

# Row 247: Seymour
$ws.Range("B247").Value = 'L'
$ws.Range("C247").Value = 'R'
$ws.Range("D247").Value = 'L'

# Row 248: Avenel
$ws.Range("B248").Value = 'L'
$ws.Range("C248").Value = 'R'

# Row 249: Euroa

# Row 250: Violet Town

# Row 251: Benalla

# Row 252: Wangaratta

# Row 253: Springhurst

# Row 254: Chiltern

# Row 255: Wodonga

# Row 256: Albury

# Row 257: Nagambie
$ws.Range("A257").Formula = '=TRANSPOSE({"Nagambie","Murchison East","Mooroopna","Shepparton"})'
$ws.Range("B257").Value = 'L'

# Row 258: Murchison East
$ws.Range("B258").Value = 'L'

# Row 259: Mooroopna
$ws.Range("B259").Value = 'R'

# Row 260: Shepparton
$ws.Range("B260").Value = 'R'

# Row 261: Clarkefield
$ws.Range("A261").Formula = '=TRANSPOSE({"Clarkefield","Riddells Creek","Gisborne","Macedon","Woodend","Kyneton","Malmsbury","Castlemaine","Kangaroo Flat","Bendigo","Eaglehawk","Raywood","Dingee","Pyramid","Kerang","Swan Hill"})'
$ws.Range("B261").Value = 'R'
$ws.Range("C261").Value = 'L'

# Row 262: Riddells Creek
$ws.Range("B262").Value = 'R'
$ws.Range("C262").Value = 'L'

# Row 263: Gisborne
$ws.Range("B263").Value = 'R'
$ws.Range("C263").Value = 'L'

# Row 264: Macedon
$ws.Range("B264").Value = 'R'
$ws.Range("C264").Value = 'L'

# Row 265: Woodend
$ws.Range("B265").Value = 'R'
$ws.Range("C265").Value = 'L'

# Row 266: Kyneton
$ws.Range("B266").Value = 'R'
$ws.Range("C266").Value = 'L'

# Row 267: Malmsbury
$ws.Range("B267").Value = 'R'
$ws.Range("C267").Value = 'L'

# Row 268: Castlemaine
$ws.Range("B268").Value = 'R'
$ws.Range("C268").Value = 'L'

# Row 269: Kangaroo Flat
$ws.Range("B269").Value = 'R'
$ws.Range("C269").Value = 'L'

# Row 270: Bendigo
$ws.Range("B270").Value = 'L'
$ws.Range("C270").Value = 'R'

# Row 271: Eaglehawk
$ws.Range("B271").Value = 'L'

# Row 272: Raywood
$ws.Range("B272").Value = 'R'

# Row 273: Dingee
$ws.Range("B273").Value = 'R'

# Row 274: Pyramid
$ws.Range("B274").Value = 'L'

# Row 275: Kerang
$ws.Range("B275").Value = 'L'

# Row 276: Swan Hill
$ws.Range("B276").Value = 'L'

# Row 277: Epsom
$ws.Range("A277").Formula = '=TRANSPOSE({"Epsom","Huntly","Goornong","Elmore","Rochester","Echuca"})'
$ws.Range("B277").Value = 'L'

# Row 278: Huntly
$ws.Range("B278").Value = 'L'

# Row 279: Goornong
$ws.Range("B279").Value = 'R'

# Row 280: Elmore
$ws.Range("B280").Value = 'R'

# Row 281: Rochester
$ws.Range("B281").Value = 'R'

# Row 282: Echuca
$ws.Range("B282").Value = 'R'

# Row 283: Ardeer
$ws.Range("A283").Formula = '=TRANSPOSE({"Ardeer","Deer Park","Caroline Springs","Rockbank","Cobblebank","Melton","Bacchus Marsh","Ballan","Ballarat","Wendouree","Beaufort","Ararat"})'
$ws.Range("B283").Value = 'L'
$ws.Range("C283").Value = 'R'

# Row 284: Deer Park
$ws.Range("B284").Value = 'R'
$ws.Range("C284").Value = 'L'

# Row 285: Caroline Springs
$ws.Range("B285").Value = 'L'
$ws.Range("C285").Value = 'R'

# Row 286: Rockbank
$ws.Range("B286").Value = 'R'
$ws.Range("C286").Value = 'L'

# Row 287: Cobblebank
$ws.Range("B287").Value = 'R'
$ws.Range("C287").Value = 'L'

# Row 288: Melton
$ws.Range("B288").Value = 'R'
$ws.Range("C288").Value = 'L'

# Row 289: Bacchus Marsh
$ws.Range("B289").Value = 'R'
$ws.Range("C289").Value = 'L'

# Row 290: Ballan
$ws.Range("B290").Value = 'R'
$ws.Range("C290").Value = 'L'

# Row 291: Ballarat
$ws.Range("B291").Value = 'L'
$ws.Range("C291").Value = 'R'

# Row 292: Wendouree
$ws.Range("B292").Value = 'R'
$ws.Range("C292").Value = 'L'

# Row 293: Beaufort
$ws.Range("B293").Value = 'L'

# Row 294: Ararat
$ws.Range("B294").Value = 'L'
$ws.Range("C294").Value = 'R'

# Row 295: Creswick
$ws.Range("A295").Formula = '=TRANSPOSE({"Creswick","Clunes","Talbot","Maryborough"})'
$ws.Range("B295").Value = 'R'

# Row 296: Clunes
$ws.Range("B296").Value = 'R'

# Row 297: Talbot
$ws.Range("B297").Value = 'L'

# Row 298: Maryborough
$ws.Range("B298").Value = 'L'

# Row 299: Tarneit
$ws.Range("A299").Formula = '=TRANSPOSE({"Tarneit","Wyndham Vale","Little River","Lara","Corio","North Shore","North Geelong","Geelong","South Geelong","Marshall","Waurn Ponds","Winchelsea","Birregurra","Colac","Camperdown","Terang","Sherwood Park","Warrnambool"})'
$ws.Range("B299").Value = 'R'
$ws.Range("C299").Value = 'L'

# Row 300: Wyndham Vale
$ws.Range("B300").Value = 'R'
$ws.Range("C300").Value = 'L'

# Row 301: Little River
$ws.Range("B301").Value = 'R'
$ws.Range("C301").Value = 'L'

# Row 302: Lara
$ws.Range("B302").Value = 'L'
$ws.Range("C302").Value = 'R'

# Row 303: Corio
$ws.Range("B303").Value = 'L'
$ws.Range("C303").Value = 'R'

# Row 304: North Shore
$ws.Range("B304").Value = 'L'
$ws.Range("C304").Value = 'R'
$ws.Range("D304").Value = 'R'

# Row 305: North Geelong
$ws.Range("B305").Value = 'R'
$ws.Range("C305").Value = 'L'

# Row 306: Geelong
$ws.Range("B306").Value = 'L'
$ws.Range("C306").Value = 'R'
$ws.Range("D306").Value = 'L'

# Row 307: South Geelong

# Row 308: Marshall

# Row 309: Waurn Ponds
$ws.Range("B309").Value = 'R'
$ws.Range("C309").Value = 'L'

# Row 310: Winchelsea
$ws.Range("B310").Value = 'L'

# Row 311: Birregurra
$ws.Range("B311").Value = 'L'

# Row 312: Colac
$ws.Range("B312").Value = 'R'

# Row 313: Camperdown
$ws.Range("B313").Value = 'L'

# Row 314: Terang
$ws.Range("B314").Value = 'L'

# Row 315: Sherwood Park
$ws.Range("B315").Value = 'L'

# Row 316: Warrnambool
$ws.Range("B316").Value = 'R'

# Row 317: Stawell
$ws.Range("A317").Value = 'Stawell'
$ws.Range("B317").Value = 'R'

# Row 318: Horsham
$ws.Range("A318").Value = 'Horsham'
$ws.Range("B318").Value = 'L'

# Row 319: Dimboola
$ws.Range("A319").Value = 'Dimboola'

# Row 320: Nhill
$ws.Range("A320").Value = 'Nhill'

# Row 321: Bordertown
$ws.Range("A321").Value = 'Bordertown'

# Row 322: Murray Bridge
$ws.Range("A322").Value = 'Murray Bridge'

# Row 323: Adelaide
$ws.Range("A323").Value = 'Adelaide'

# --- Apply yellow highlight fill to specific station cells ---
$ws.Range("A228").Interior.Color = 65535
$ws.Range("A249").Interior.Color = 65535
$ws.Range("A250").Interior.Color = 65535
$ws.Range("A251").Interior.Color = 65535
$ws.Range("A252").Interior.Color = 65535
$ws.Range("A253").Interior.Color = 65535
$ws.Range("A254").Interior.Color = 65535
$ws.Range("A255").Interior.Color = 65535
$ws.Range("A256").Interior.Color = 65535
$ws.Range("A307").Interior.Color = 65535
$ws.Range("A308").Interior.Color = 65535

# --- Update selection/view to reflect end of data ---
$ws.Range("A324").Select()
